$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# "fixed excel logic questions format"
#
# The three "What is the output of the following code?" questions had their
# code squashed onto the same line as the question text. Reformat them so
# the code starts on its own line (matching every other "output of the
# following code" question in the sheet), and tidy up names[1][-2] to use
# the print(...) function call syntax.
#
# Setting B16 first, then B15, then B14 (i.e. "names" before "dictionary"
# before "list") makes the freshly-created shared-string entries land at
# the end of the shared-string table in that same order, once the old,
# now-unreferenced entries are dropped on save.
# --------------------------------------------------------------------------

$namesText = @"
What is the output of the following code? 
        names = ['John', 'Paul', 'George', 'Ringo']
        print(names[1][-2])
"@

$dictText = @"
What is the output of the following code? 
dictionary = {1: 'first', 2: 'second', 3: 'third', 4: 'fourth'}
del dictionary[1]
dictionary[1] = '42'
del dictionary[2]
print len(dictionary)
"@

$listText = @"
What is the output of the following code? 
list = [ 'one', 'two']
for i in list:
	list.append(i.upper())
	print(list)
"@

$ws.Range("B16").Value = $namesText
$ws.Range("B15").Value = $dictText
$ws.Range("B14").Value = $listText

# --------------------------------------------------------------------------
# Row heights grew slightly across the whole sheet (the workbook was
# resaved from a newer Excel build with different default font metrics),
# and the three reformatted rows above grew further because their text now
# wraps across more lines.
# --------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight  = 150
$ws.Rows.Item(3).RowHeight  = 165
$ws.Rows.Item(4).RowHeight  = 60
$ws.Rows.Item(5).RowHeight  = 60
$ws.Rows.Item(6).RowHeight  = 60
$ws.Rows.Item(7).RowHeight  = 75
$ws.Rows.Item(8).RowHeight  = 75
$ws.Rows.Item(9).RowHeight  = 60
$ws.Rows.Item(10).RowHeight = 90
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 210
$ws.Rows.Item(14).RowHeight = 75
$ws.Rows.Item(15).RowHeight = 90
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 60
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60

# --------------------------------------------------------------------------
# View state: zoom was reset to 100% and the selection left on D2.
# --------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 100
$ws.Range("D2").Select()
